# "json to excel, reinicio de entradas post registro"
# Append newly-registered income/expense entries to the bottom of each
# sheet's data table (post "form" submission), then reset the entry rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Ingresos": two new income entries registered on 2025-05-06
# ---------------------------------------------------------------------
$wsIngresos = $wb.Worksheets.Item("Ingresos")

$wsIngresos.Range("A4").Value = "Ingreso"
$wsIngresos.Range("B4").Value = "Otros"
$wsIngresos.Range("C4").Value = "transferencia"
$wsIngresos.Range("D4").Value = 2675

$wsIngresos.Range("A5").Value = "Ingreso"
$wsIngresos.Range("B5").Value = "Inversiones"
$wsIngresos.Range("C5").Value = "nu 14% 3 meses"
$wsIngresos.Range("D5").Value = 766

# "fecha" is stored as plain text, not a date serial, everywhere else in
# the sheet. Force a Text number format before assigning so Excel does
# not auto-convert the date-looking string, then strip the formatting
# back off so the cells end up with no explicit style, matching the
# rest of the un-styled data rows.
$rngFechaIngresos = $wsIngresos.Range("E4:E5")
$rngFechaIngresos.NumberFormat = "@"
$wsIngresos.Range("E4").Value = "2025-05-06"
$wsIngresos.Range("E5").Value = "2025-05-06"
$rngFechaIngresos.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "Egresos": three new expense entries registered on 2025-05-06
# ---------------------------------------------------------------------
$wsEgresos = $wb.Worksheets.Item("Egresos")

$wsEgresos.Range("A4").Value = "Egreso"
$wsEgresos.Range("B4").Value = "Vivienda"
$wsEgresos.Range("C4").Value = "servicio de gas"
$wsEgresos.Range("D4").Value = 550

$wsEgresos.Range("A5").Value = "Egreso"
$wsEgresos.Range("B5").Value = "Salud"
$wsEgresos.Range("C5").Value = "medicamentos"
$wsEgresos.Range("D5").Value = 1250

$wsEgresos.Range("A6").Value = "Egreso"
$wsEgresos.Range("B6").Value = "Alimentos"
$wsEgresos.Range("C6").Value = "café"
$wsEgresos.Range("D6").Value = 360

$rngFechaEgresos = $wsEgresos.Range("E4:E6")
$rngFechaEgresos.NumberFormat = "@"
$wsEgresos.Range("E4").Value = "2025-05-06"
$wsEgresos.Range("E5").Value = "2025-05-06"
$wsEgresos.Range("E6").Value = "2025-05-06"
$rngFechaEgresos.ClearFormats()
